# stateless entities outside the US
#
# The "IMF (20%)" model columns are inserted right after the "GFA" columns
# (and right before the existing "IMF" columns) in each of the three
# 8-column model blocks (B:I, J:Q, R:Y). This pushes the old "OECD (20%)"
# columns out entirely (dropped) and the former "IMF" columns shift two
# slots to the right, replacing the old "OECD (20%)" header/position.
#
# Net effect per 8-column block:
#   col3 (old IMF-Sales)        -> becomes new "IMF (20%) - Sales" data
#   col4 (old IMF-Sales+Emp)    -> becomes new "IMF (20%) - Sales + Emp" data
#   col5 (old OECD(20%)-Sales)      -> replaced with old col3 value (IMF - Sales)
#   col6 (old OECD(20%)-Sales+Emp)  -> replaced with old col4 value (IMF - Sales + Emp)
#   col7/col8 (OECD - Sales / OECD - Sales + Emp) -> unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: header labels for each of the three 8-column model blocks
# ---------------------------------------------------------------------
$headers = @("GFA - Sales", "GFA - Sales + Emp", "IMF (20%) - Sales", "IMF (20%) - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")

$blockStarts = @(2, 10, 18)   # column numbers for B, J, R
foreach ($startCol in $blockStarts) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(2, $startCol + $i).Value = $headers[$i]
    }
}

# ---------------------------------------------------------------------
# Data rows: 4, 6, 7, 8 — new IMF(20%) values (col D/L/T, E/M/U) and the
# shifted-right former-IMF values (col F/N/V, G/O/W).
# ---------------------------------------------------------------------

# New, previously-nonexistent "IMF (20%)" raw data (block B:I only — this
# is genuinely new source data, not derivable from other cells).
$ws.Range("D4").Value = 0.010903556787023
$ws.Range("E4").Value = 0.009036458532427226

$ws.Range("D6").Value = 0.007050127574638287
$ws.Range("E6").Value = 0.02773577005388293

$ws.Range("D7").Value = 0.006942562708641375
$ws.Range("E7").Value = 0.007822019816789091

$ws.Range("D8").Value = 0.002250293725144459
$ws.Range("E8").Value = 0.005258857082920571

# Former "IMF - Sales" / "IMF - Sales + Emp" raw values shift right into
# F/G (old OECD(20%) slot), replacing the dropped OECD(20%) raw data.
$ws.Range("F4").Value = 0.05451778393511495
$ws.Range("G4").Value = 0.04518229266213623

$ws.Range("F6").Value = 0.03525063787319145
$ws.Range("G6").Value = 0.1386788502694146

$ws.Range("F7").Value = 0.03471281354320685
$ws.Range("G7").Value = 0.03911009908394546

$ws.Range("F8").Value = 0.01125146862572229
$ws.Range("G8").Value = 0.02629428541460285

# Block J:Q (M_ETR) — the IMF columns (L/M) keep their values; the new
# N/O slot (old OECD(20%) position) is simply a copy of L/M.
$ws.Range("N4").Value = $ws.Range("L4").Value()
$ws.Range("O4").Value = $ws.Range("M4").Value()

$ws.Range("N6").Value = $ws.Range("L6").Value()
$ws.Range("O6").Value = $ws.Range("M6").Value()

$ws.Range("N7").Value = $ws.Range("L7").Value()
$ws.Range("O7").Value = $ws.Range("M7").Value()

$ws.Range("N8").Value = $ws.Range("L8").Value()
$ws.Range("O8").Value = $ws.Range("M8").Value()

# Block R:Y (M_PL) — same pattern: T/U (IMF) values copied into V/W.
$ws.Range("V4").Value = $ws.Range("T4").Value()
$ws.Range("W4").Value = $ws.Range("U4").Value()

$ws.Range("V6").Value = $ws.Range("T6").Value()
$ws.Range("W6").Value = $ws.Range("U6").Value()

$ws.Range("V7").Value = $ws.Range("T7").Value()
$ws.Range("W7").Value = $ws.Range("U7").Value()

$ws.Range("V8").Value = $ws.Range("T8").Value()
$ws.Range("W8").Value = $ws.Range("U8").Value()
